$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for Wins / Losses / Ties (columns AD, AE, AF on row 1),
# matching the bold/centered/bordered style already used by the other
# header cells (copy formatting from the neighboring AC1 header cell).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for each data row (2 through 55) with the same
# values: 67 wins, 95 losses, 0 ties.
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 67   # AD
    $ws.Cells.Item($r, 31).Value = 95   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
